# Apply updated values to team matrix (Bryan_B) from games pulled March 7.
# Updates probability values in rows 2-19 (columns B through S) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 0.2414634146341463
    $ws.Range("C2").Value = 0.4878048780487805
    $ws.Range("J2").Value = 0.01707317073170732
    $ws.Range("P2").Value = 0.1731707317073171
    $ws.Range("S2").Value = 0.08048780487804878
    $ws.Range("B3").Value = 0.02252252252252252
    $ws.Range("C3").Value = 0.09009009009009009
    $ws.Range("J3").Value = 0.04954954954954955
    $ws.Range("P3").Value = 0.6936936936936937
    $ws.Range("S3").Value = 0.1441441441441441
    $ws.Range("J4").Value = 0.04081632653061224
    $ws.Range("P4").Value = 0.5510204081632653
    $ws.Range("S4").Value = 0.4081632653061225
    $ws.Range("B6").Value = 0.08071748878923767
    $ws.Range("D6").Value = 0.01345291479820628
    $ws.Range("F6").Value = 0.03139013452914798
    $ws.Range("J6").Value = 0.2376681614349776
    $ws.Range("O6").Value = 0.04035874439461883
    $ws.Range("Q6").Value = 0.1345291479820628
    $ws.Range("R6").Value = 0.07174887892376682
    $ws.Range("S6").Value = 0.3901345291479821
    $ws.Range("B7").Value = 0.102803738317757
    $ws.Range("D7").Value = 0.02336448598130841
    $ws.Range("F7").Value = 0.0514018691588785
    $ws.Range("J7").Value = 0.1401869158878505
    $ws.Range("O7").Value = 0.02336448598130841
    $ws.Range("Q7").Value = 0.1495327102803738
    $ws.Range("R7").Value = 0.07009345794392523
    $ws.Range("S7").Value = 0.4392523364485981
    $ws.Range("B8").Value = 0.06806282722513089
    $ws.Range("D8").Value = 0.02443280977312391
    $ws.Range("E8").Value = 0.001745200698080279
    $ws.Range("F8").Value = 0.04013961605584642
    $ws.Range("J8").Value = 0.1308900523560209
    $ws.Range("O8").Value = 0.02792321116928447
    $ws.Range("Q8").Value = 0.1797556719022688
    $ws.Range("R8").Value = 0.07678883071553229
    $ws.Range("S8").Value = 0.450261780104712
    $ws.Range("B9").Value = 0.08171206225680934
    $ws.Range("D9").Value = 0.01945525291828794
    $ws.Range("F9").Value = 0.05058365758754864
    $ws.Range("J9").Value = 0.1050583657587549
    $ws.Range("O9").Value = 0.01945525291828794
    $ws.Range("Q9").Value = 0.2062256809338521
    $ws.Range("R9").Value = 0.07392996108949416
    $ws.Range("S9").Value = 0.443579766536965
    $ws.Range("B10").Value = 0.1172055427251732
    $ws.Range("D10").Value = 0.01501154734411085
    $ws.Range("F10").Value = 0.05311778290993072
    $ws.Range("J10").Value = 0.1431870669745959
    $ws.Range("O10").Value = 0.0207852193995381
    $ws.Range("Q10").Value = 0.1980369515011547
    $ws.Range("R10").Value = 0.08025404157043881
    $ws.Range("S10").Value = 0.3724018475750577
    $ws.Range("G11").Value = 0.16
    $ws.Range("J11").Value = 0.12
    $ws.Range("K11").Value = 0.2314285714285714
    $ws.Range("L11").Value = 0.4742857142857143
    $ws.Range("S11").Value = 0.01428571428571429
    $ws.Range("G12").Value = 0.7237569060773481
    $ws.Range("J12").Value = 0.1823204419889503
    $ws.Range("K12").Value = 0.005524861878453038
    $ws.Range("L12").Value = 0.06629834254143646
    $ws.Range("S12").Value = 0.02209944751381215
    $ws.Range("G13").Value = 0.55
    $ws.Range("J13").Value = 0.4166666666666667
    $ws.Range("S13").Value = 0.03333333333333333
    $ws.Range("F15").Value = 0.0234375
    $ws.Range("H15").Value = 0.1640625
    $ws.Range("I15").Value = 0.05078125
    $ws.Range("J15").Value = 0.33984375
    $ws.Range("K15").Value = 0.03125
    $ws.Range("M15").Value = 0.015625
    $ws.Range("O15").Value = 0.06640625
    $ws.Range("S15").Value = 0.30859375
    $ws.Range("F16").Value = 0.02904564315352697
    $ws.Range("H16").Value = 0.1701244813278008
    $ws.Range("I16").Value = 0.07468879668049792
    $ws.Range("J16").Value = 0.4232365145228216
    $ws.Range("K16").Value = 0.1037344398340249
    $ws.Range("M16").Value = 0.01244813278008299
    $ws.Range("O16").Value = 0.03734439834024896
    $ws.Range("S16").Value = 0.1493775933609958
    $ws.Range("F17").Value = 0.01423487544483986
    $ws.Range("H17").Value = 0.1832740213523132
    $ws.Range("I17").Value = 0.1085409252669039
    $ws.Range("J17").Value = 0.4644128113879004
    $ws.Range("K17").Value = 0.08896797153024912
    $ws.Range("M17").Value = 0.01779359430604982
    $ws.Range("N17").Value = 0.001779359430604982
    $ws.Range("O17").Value = 0.05160142348754448
    $ws.Range("S17").Value = 0.0693950177935943
    $ws.Range("F18").Value = 0.004347826086956522
    $ws.Range("H18").Value = 0.2260869565217391
    $ws.Range("I18").Value = 0.06086956521739131
    $ws.Range("J18").Value = 0.4434782608695652
    $ws.Range("K18").Value = 0.1043478260869565
    $ws.Range("M18").Value = 0.01739130434782609
    $ws.Range("O18").Value = 0.08695652173913043
    $ws.Range("S18").Value = 0.05652173913043478
    $ws.Range("F19").Value = 0.01201011378002529
    $ws.Range("H19").Value = 0.2149178255372946
    $ws.Range("I19").Value = 0.09671302149178256
    $ws.Range("J19").Value = 0.411504424778761
    $ws.Range("K19").Value = 0.1005056890012642
    $ws.Range("M19").Value = 0.02528445006321112
    $ws.Range("N19").Value = 0.0006321112515802782
    $ws.Range("O19").Value = 0.05246523388116309
    $ws.Range("S19").Value = 0.08596713021491782

$wb.Save()
